# Update the "Urbanization" breakdown labels (rows 23-24) on the active sheet
# from the short "City/Village" style names to the fuller
# "Urban settlements / Rural area" style names, in all three languages
# (Kyrgyz, Russian, English columns A/B/C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: "Шаар" / "Город" / "Urban"  ->  "Шаар жерлери" / "Городские поселения" / "City"
$ws.Range("A23").Value = "Шаар жерлери"
$ws.Range("B23").Value = "Городские поселения"
$ws.Range("C23").Value = "City"

# Row 24: "Айыл" / "Село" / "Rural"  ->  "Айыл аймагы" / "Сельская местность" / "Village"
$ws.Range("A24").Value = "Айыл аймагы"
$ws.Range("B24").Value = "Сельская местность"
$ws.Range("C24").Value = "Village"

# Move the active selection to C30, matching the refreshed view state.
$ws.Range("C30").Select()
